$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: change style from (4,5) to (8,9) by copying formats from row 11,
#     which already uses style 8 (A/B) and style 9 (C/D/E). This reuses the
#     existing cellXf entries instead of minting new style entries.
$ws.Range("A11:B11").Copy() | Out-Null
$ws.Range("A14:B14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C11:E11").Copy() | Out-Null
$ws.Range("C14:E14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- New rows 15-16: fill column by column (filename, then English col,
#     then Russian col, then converted/encoded col; numbers can go anytime)
#     so the shared-string table is built up in the same order as the
#     authored workbook.
$ws.Range("A15").Value = 'SCRIPT/T01P02A/us2308.ssb'
$ws.Range("B15").Value = 19
$ws.Range("B16").Value = 22

$ws.Range("C15").Value = ' Exploring is far too much\neffort...[K] Maybe we\''ll go to the Hot Spring…'
$ws.Range("C16").Value = ' But going to the Hot Spring from\nhere takes too much effort…'

$ws.Range("D15").Value = ' Исследования утомляют...[K] Может,\nсходить к Горячим Источникам...'
$ws.Range("D16").Value = ' Но путь к Горячим Источникам\nтак утомителен...'

$ws.Range("E15").Value = ' Éòòìåäïâàîéÿ ôóïíìÿýó...[K] Íïçåó,\nòöïäéóû ë Ãïñÿœéí Éòóïœîéëàí...'
$ws.Range("E16").Value = ' Îï ðôóû ë Ãïñÿœéí Éòóïœîéëàí\nóàë ôóïíéóåìåî...'

# --- Row heights to mirror the authored layout ---
$ws.Rows.Item(15).RowHeight = 43.2
$ws.Rows.Item(16).RowHeight = 21.6

# --- Selection / scroll position ---
$ws.Activate()
$ws.Range("D14").Select()
